$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to re-pulled data per commit message:
# "repull data, push all data, mean calculation"
$ws.Range("F6").Value = -4
$ws.Range("F9").Value = 0
$ws.Range("F17").Value = -4
$ws.Range("F30").Value = -2
$ws.Range("F41").Value = -2
$ws.Range("F45").Value = -6
$ws.Range("F52").Value = -5
$ws.Range("F55").Value = 6
$ws.Range("F59").Value = -1
$ws.Range("F62").Value = -2
$ws.Range("F64").Value = 11
